# Update column G ("K") values for rows 2-34 on the active worksheet.
# These values were regenerated to use K (strikeouts) instead of Strike#,
# replacing the previous raw "Strike#" figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 4
    3  = 2
    4  = 2
    5  = 6
    6  = 3
    7  = 8
    8  = 4
    9  = 5
    10 = 10
    11 = 9
    12 = 4
    13 = 3
    14 = 4
    15 = 3
    16 = 2
    17 = 4
    18 = 5
    19 = 10
    20 = 7
    21 = 3
    22 = 4
    23 = 6
    24 = 4
    25 = 5
    26 = 5
    27 = 4
    28 = 9
    29 = 4
    30 = 3
    31 = 4
    32 = 5
    33 = 5
    34 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
